# Update Fresno St._B team-specific transition matrix with new probabilities
# (more games were simulated, shifting the per-state outcome distributions).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2412280701754386
$ws.Range("C2").Value = 0.4649122807017544
$ws.Range("J2").Value = 0.0131578947368421
$ws.Range("P2").Value = 0.1798245614035088
$ws.Range("S2").Value = 0.1008771929824561

# Row 3
$ws.Range("B3").Value = 0.009259259259259259
$ws.Range("C3").Value = 0.01851851851851852
$ws.Range("P3").Value = 0.7222222222222222
$ws.Range("S3").Value = 0.25

# Row 4
$ws.Range("J4").Value = 0.06666666666666667
$ws.Range("P4").Value = 0.6444444444444445
$ws.Range("S4").Value = 0.2888888888888889

# Row 6
$ws.Range("B6").Value = 0.0797872340425532
$ws.Range("D6").Value = 0.01063829787234043
$ws.Range("F6").Value = 0.03191489361702127
$ws.Range("J6").Value = 0.3031914893617021
$ws.Range("O6").Value = 0.01595744680851064
$ws.Range("Q6").Value = 0.1861702127659574
$ws.Range("R6").Value = 0.05319148936170213
$ws.Range("S6").Value = 0.3191489361702128

# Row 7
$ws.Range("B7").Value = 0.1027397260273973
$ws.Range("D7").Value = 0.0136986301369863
$ws.Range("E7").Value = 0.00684931506849315
$ws.Range("F7").Value = 0.07534246575342465
$ws.Range("J7").Value = 0.0958904109589041
$ws.Range("O7").Value = 0.02054794520547945
$ws.Range("Q7").Value = 0.1780821917808219
$ws.Range("R7").Value = 0.1232876712328767
$ws.Range("S7").Value = 0.3835616438356164

# Row 8
$ws.Range("B8").Value = 0.05263157894736842
$ws.Range("D8").Value = 0.01503759398496241
$ws.Range("E8").Value = 0.007518796992481203
$ws.Range("F8").Value = 0.05513784461152882
$ws.Range("J8").Value = 0.1353383458646616
$ws.Range("O8").Value = 0.03007518796992481
$ws.Range("Q8").Value = 0.2005012531328321
$ws.Range("R8").Value = 0.112781954887218
$ws.Range("S8").Value = 0.3909774436090225

# Row 9
$ws.Range("B9").Value = 0.07602339181286549
$ws.Range("D9").Value = 0.02339181286549707
$ws.Range("F9").Value = 0.08771929824561403
$ws.Range("J9").Value = 0.1052631578947368
$ws.Range("O9").Value = 0.01169590643274854
$ws.Range("Q9").Value = 0.1929824561403509
$ws.Range("R9").Value = 0.1345029239766082
$ws.Range("S9").Value = 0.3684210526315789

# Row 10
$ws.Range("B10").Value = 0.09502664298401421
$ws.Range("D10").Value = 0.02841918294849023
$ws.Range("F10").Value = 0.05772646536412078
$ws.Range("J10").Value = 0.127886323268206
$ws.Range("O10").Value = 0.01154529307282416
$ws.Range("Q10").Value = 0.2211367673179396
$ws.Range("R10").Value = 0.08081705150976909
$ws.Range("S10").Value = 0.3774422735346359

# Row 11
$ws.Range("G11").Value = 0.1506276150627615
$ws.Range("J11").Value = 0.1129707112970711
$ws.Range("K11").Value = 0.2217573221757322
$ws.Range("L11").Value = 0.502092050209205
$ws.Range("S11").Value = 0.01255230125523013

# Row 12
$ws.Range("G12").Value = 0.7109375
$ws.Range("J12").Value = 0.203125
$ws.Range("K12").Value = 0.015625
$ws.Range("L12").Value = 0.0390625
$ws.Range("S12").Value = 0.03125

# Row 13
$ws.Range("G13").Value = 0.5454545454545454
$ws.Range("J13").Value = 0.4545454545454545

# Row 14
$ws.Range("G14").Value = 0.75
$ws.Range("J14").Value = 0.25

# Row 15
$ws.Range("F15").Value = 0.02
$ws.Range("H15").Value = 0.155
$ws.Range("I15").Value = 0.08
$ws.Range("J15").Value = 0.42
$ws.Range("K15").Value = 0.05
$ws.Range("M15").Value = 0.015
$ws.Range("N15").Value = 0.005
$ws.Range("S15").Value = 0.205

# Row 16
$ws.Range("F16").Value = 0.02040816326530612
$ws.Range("H16").Value = 0.1904761904761905
$ws.Range("I16").Value = 0.1428571428571428
$ws.Range("J16").Value = 0.3673469387755102
$ws.Range("K16").Value = 0.07482993197278912
$ws.Range("M16").Value = 0.01360544217687075
$ws.Range("O16").Value = 0.06802721088435375
$ws.Range("S16").Value = 0.1224489795918367

# Row 17
$ws.Range("F17").Value = 0.02612826603325416
$ws.Range("H17").Value = 0.1995249406175772
$ws.Range("I17").Value = 0.07363420427553444
$ws.Range("J17").Value = 0.4418052256532066
$ws.Range("K17").Value = 0.08788598574821853
$ws.Range("M17").Value = 0.01900237529691211
$ws.Range("N17").Value = 0.002375296912114014
$ws.Range("O17").Value = 0.05700712589073634
$ws.Range("S17").Value = 0.09263657957244656

# Row 18
$ws.Range("F18").Value = 0.02150537634408602
$ws.Range("H18").Value = 0.1935483870967742
$ws.Range("I18").Value = 0.07526881720430108
$ws.Range("J18").Value = 0.3709677419354839
$ws.Range("K18").Value = 0.1236559139784946
$ws.Range("M18").Value = 0.01612903225806452
$ws.Range("N18").Value = 0.005376344086021506
$ws.Range("O18").Value = 0.07526881720430108
$ws.Range("S18").Value = 0.1182795698924731

# Row 19
$ws.Range("F19").Value = 0.02455795677799607
$ws.Range("H19").Value = 0.2170923379174853
$ws.Range("I19").Value = 0.08742632612966601
$ws.Range("J19").Value = 0.37524557956778
$ws.Range("K19").Value = 0.1011787819253438
$ws.Range("M19").Value = 0.02652259332023576
$ws.Range("N19").Value = 0.001964636542239686
$ws.Range("O19").Value = 0.08153241650294696
$ws.Range("S19").Value = 0.08447937131630648
